$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 15 (Configuration Item #14 - "Configuration Item List")
#   C15: 1.2 -> 1.3
#   D15: 41927 -> 41928  (2014-10-15 -> 2014-10-16)
#   F15: "ConfigurationItemList.xlsx" -> hyperlink to full github URL
# ---------------------------------------------------------------
$ws.Range("C15").Value = 1.3
$ws.Range("D15").Value = 41928

$ws.Hyperlinks.Add(
    $ws.Range("F15"),
    "https://github.com/sungori/SoftBugOff/blob/master/ConfigurationItemList.xlsx",
    [Type]::Missing,
    [Type]::Missing,
    "https://github.com/sungori/SoftBugOff/blob/master/ConfigurationItemList.xlsx"
) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F15").PasteSpecial(-4122) | Out-Null

$ws.Rows.Item(15).RowHeight = 30.75

# ---------------------------------------------------------------
# Row 14 (Configuration Item #13 - "Term Project Document")
#   C14: 0.5 -> 0.6
#   F14: empty -> hyperlink to Giles_Vernon_TermProject.docx
#   Row height: 15.75 -> 30.75
# ---------------------------------------------------------------
$ws.Range("C14").Value = 0.6

$ws.Hyperlinks.Add(
    $ws.Range("F14"),
    "https://github.com/sungori/SoftBugOff/blob/master/Giles_Vernon_TermProject.docx",
    [Type]::Missing,
    [Type]::Missing,
    "https://github.com/sungori/SoftBugOff/blob/master/Giles_Vernon_TermProject.docx"
) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F14").PasteSpecial(-4122) | Out-Null

$ws.Rows.Item(14).RowHeight = 30.75

# ---------------------------------------------------------------
# Row 8 (Configuration Item #7 - "Definition of use cases")
#   C8: 0.2 -> 0.3
#   D8: 41921 -> 41928  (2014-10-09 -> 2014-10-16)
#   F8: "This Document (Section 2.2)" -> hyperlink to UseCases.xlsx
#   Row height: 15.75 -> 30.75
# ---------------------------------------------------------------
$ws.Range("C8").Value = 0.3
$ws.Range("D8").Value = 41928

$ws.Hyperlinks.Add(
    $ws.Range("F8"),
    "https://github.com/sungori/SoftBugOff/blob/master/UseCases.xlsx",
    [Type]::Missing,
    [Type]::Missing,
    "https://github.com/sungori/SoftBugOff/blob/master/UseCases.xlsx"
) | Out-Null
# Restore the original local "hyperlink-look" cell style (border/wrap/etc.)
# that Hyperlinks.Add overwrote, matching the style used by the other
# hyperlinked cells in column F (e.g. F2).
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F8").PasteSpecial(-4122) | Out-Null

$ws.Rows.Item(8).RowHeight = 30.75

# ---------------------------------------------------------------
# View state: scroll so row 4 is at the top, select C8
# ---------------------------------------------------------------
$ws.Range("C8").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
